# Arizona overview workbook restructuring
# - Reorders the metric columns so "Share of 990 filers with government
#   grants at risk" becomes the first metric column (right after any
#   geography/size/subsector label column).
# - Renames "Operating surplus with/without government grants (%)" to
#   "Size of operating surplus with/without government grants".
# - Renames several row labels (Congressional districts, size buckets,
#   subsectors) and reorders a couple of rows (Size sheet, Subsector sheet).

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $row, $col, $value) {
    # Force the cell to stay plain text even when the string looks like a
    # number/currency/percentage (e.g. "1,446", "$3,710,995,891", "65.98%"),
    # matching the source workbook's inlineStr cells.
    $ws.Cells.Item($row, $col).NumberFormat = "@"
    $ws.Cells.Item($row, $col).Value = $value
}

function Set-LabelValue($ws, $row, $col, $value) {
    # Plain (non-numeric-looking) text labels/headers - no NumberFormat
    # fiddling needed, keeps the cell's existing style untouched.
    $ws.Cells.Item($row, $col).Value = $value
}

function Get-Text($ws, $row, $col) {
    return $ws.Cells.Item($row, $col).Text
}

# ---------------------------------------------------------------------------
# Sheet "Overall": 5 metric columns, no leading label column.
#   OLD: A=Number B=Total$ C=OpWith% D=OpWithout% E=ShareAtRisk
#   NEW: A=ShareAtRisk B=Number C=Total$ D=SizeOpWith E=SizeOpWithout
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overall")

Set-LabelValue $ws 1 1 "Share of 990 filers with government grants at risk"
Set-LabelValue $ws 1 2 "Number of 990 filers with government grants"
Set-LabelValue $ws 1 3 "Total government grants (`$)"
Set-LabelValue $ws 1 4 "Size of operating surplus with government grants"
Set-LabelValue $ws 1 5 "Size of operating surplus without government grants"

$oldA = Get-Text $ws 2 1
$oldB = Get-Text $ws 2 2
$oldC = Get-Text $ws 2 3
$oldD = Get-Text $ws 2 4
$oldE = Get-Text $ws 2 5

Set-TextValue $ws 2 1 $oldE
Set-TextValue $ws 2 2 $oldA
Set-TextValue $ws 2 3 $oldB
Set-TextValue $ws 2 4 $oldC
Set-TextValue $ws 2 5 $oldD

# ---------------------------------------------------------------------------
# Sheets with a leading label column (Geography / Size / Subsector) and
# five metric columns B..F:
#   OLD: B=Number C=Total$ D=OpWith% E=OpWithout% F=ShareAtRisk
#   NEW: B=ShareAtRisk C=Number D=Total$ E=OpWith% F=OpWithout%
#        (header row also renames D/E -> "Size of operating surplus ...")
# ---------------------------------------------------------------------------
function Update-MetricHeader($ws) {
    Set-LabelValue $ws 1 2 "Share of 990 filers with government grants at risk"
    Set-LabelValue $ws 1 3 "Number of 990 filers with government grants"
    Set-LabelValue $ws 1 4 "Total government grants (`$)"
    Set-LabelValue $ws 1 5 "Size of operating surplus with government grants"
    Set-LabelValue $ws 1 6 "Size of operating surplus without government grants"
}

function Shift-MetricRow($ws, $row) {
    $oldB = Get-Text $ws $row 2
    $oldC = Get-Text $ws $row 3
    $oldD = Get-Text $ws $row 4
    $oldE = Get-Text $ws $row 5
    $oldF = Get-Text $ws $row 6

    Set-TextValue $ws $row 2 $oldF
    Set-TextValue $ws $row 3 $oldB
    Set-TextValue $ws $row 4 $oldC
    Set-TextValue $ws $row 5 $oldD
    Set-TextValue $ws $row 6 $oldE
}

# ---------------------------------------------------------------------------
# Sheet "County": label column (A) text is unchanged; only the metric
# columns shift.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("County")
Update-MetricHeader $ws
for ($r = 2; $r -le 17; $r++) {
    Shift-MetricRow $ws $r
}

# ---------------------------------------------------------------------------
# Sheet "Congressional District": metric columns shift, AND the district
# labels change from "Nth Congressional district" to
# "Congressional District N" (rows 4-12; rows 2-3 are US/Arizona totals).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Congressional District")
Update-MetricHeader $ws
for ($r = 2; $r -le 12; $r++) {
    Shift-MetricRow $ws $r
}
Set-LabelValue $ws 4  1 "Congressional District 1"
Set-LabelValue $ws 5  1 "Congressional District 2"
Set-LabelValue $ws 6  1 "Congressional District 3"
Set-LabelValue $ws 7  1 "Congressional District 4"
Set-LabelValue $ws 8  1 "Congressional District 5"
Set-LabelValue $ws 9  1 "Congressional District 6"
Set-LabelValue $ws 10 1 "Congressional District 7"
Set-LabelValue $ws 11 1 "Congressional District 8"
Set-LabelValue $ws 12 1 "Congressional District 9"

# ---------------------------------------------------------------------------
# Sheet "Size": metric columns shift, labels are renamed, AND rows are
# reordered into ascending size order (rows 2-7; row 8 is the Total).
#   New row order (label -> which OLD row it came from):
#     2 Between $100K and $499K   <- old row 2 ($100K to $499K)
#     3 Between $1M and $4.99M    <- old row 4 ($1M to $4.9M)
#     4 Between $500K and $999K   <- old row 5 ($500K to $999K)
#     5 Between $5M and $9.99M    <- old row 6 ($5M to $9.9M)
#     6 Greater than $10M         <- old row 3 ($10M or more)
#     7 Less than $100K           <- old row 7 (unchanged label)
#     8 Total                     <- old row 8 (unchanged)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Size")
Update-MetricHeader $ws

# Snapshot every data row (label + the 5 metric columns, still in OLD
# column order B..F) before any writes, so the reorder can't clobber
# values it still needs to read.
$sizeRows = @{}
for ($r = 2; $r -le 8; $r++) {
    $sizeRows[$r] = @(
        (Get-Text $ws $r 1),
        (Get-Text $ws $r 2),
        (Get-Text $ws $r 3),
        (Get-Text $ws $r 4),
        (Get-Text $ws $r 5),
        (Get-Text $ws $r 6)
    )
}

function Write-SizeRow($ws, $destRow, $label, $srcRowData) {
    # $srcRowData = @(label, B, C, D, E, F) in the OLD column layout.
    Set-LabelValue $ws $destRow 1 $label
    Set-TextValue $ws $destRow 2 $srcRowData[5]   # F (ShareAtRisk) -> B
    Set-TextValue $ws $destRow 3 $srcRowData[1]   # B (Number)      -> C
    Set-TextValue $ws $destRow 4 $srcRowData[2]   # C (Total$)      -> D
    Set-TextValue $ws $destRow 5 $srcRowData[3]   # D (OpWith%)     -> E
    Set-TextValue $ws $destRow 6 $srcRowData[4]   # E (OpWithout%)  -> F
}

Write-SizeRow $ws 2 "Between `$100K and `$499K"   $sizeRows[2]
Write-SizeRow $ws 3 "Between `$1M and `$4.99M"    $sizeRows[4]
Write-SizeRow $ws 4 "Between `$500K and `$999K"   $sizeRows[5]
Write-SizeRow $ws 5 "Between `$5M and `$9.99M"    $sizeRows[6]
Write-SizeRow $ws 6 "Greater than `$10M"          $sizeRows[3]
Write-SizeRow $ws 7 "Less than `$100K"            $sizeRows[7]
Write-SizeRow $ws 8 "Total"                       $sizeRows[8]

# ---------------------------------------------------------------------------
# Sheet "Subsector": metric columns shift, several labels are renamed, and
# the Universities/Unclassified rows swap places (rows 11 <-> 12; row 13
# is the Total).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Subsector")
Update-MetricHeader $ws
for ($r = 2; $r -le 13; $r++) {
    Shift-MetricRow $ws $r
}

# Rename labels that changed capitalization / wording (rows stay put).
Set-LabelValue $ws 2  1 "Arts, Culture, and Humanities"
Set-LabelValue $ws 3  1 "Education (Excluding Universities)"
Set-LabelValue $ws 4  1 "Environment and Animals"
Set-LabelValue $ws 5  1 "Health (Excluding Hospitals)"
Set-LabelValue $ws 7  1 "Human Services"
Set-LabelValue $ws 8  1 "International, Foreign Affairs"
Set-LabelValue $ws 9  1 "Public, Societal Benefit"
Set-LabelValue $ws 10 1 "Religion Related"

# Swap the Universities (old row 11) and Unclassified (old row 12) rows -
# snapshot both (now in the NEW column orderB..F) before overwriting.
$row11 = @(
    (Get-Text $ws 11 1), (Get-Text $ws 11 2), (Get-Text $ws 11 3),
    (Get-Text $ws 11 4), (Get-Text $ws 11 5), (Get-Text $ws 11 6)
)
$row12 = @(
    (Get-Text $ws 12 1), (Get-Text $ws 12 2), (Get-Text $ws 12 3),
    (Get-Text $ws 12 4), (Get-Text $ws 12 5), (Get-Text $ws 12 6)
)

Set-LabelValue $ws 11 1 $row12[0]
Set-TextValue  $ws 11 2 $row12[1]
Set-TextValue  $ws 11 3 $row12[2]
Set-TextValue  $ws 11 4 $row12[3]
Set-TextValue  $ws 11 5 $row12[4]
Set-TextValue  $ws 11 6 $row12[5]

Set-LabelValue $ws 12 1 $row11[0]
Set-TextValue  $ws 12 2 $row11[1]
Set-TextValue  $ws 12 3 $row11[2]
Set-TextValue  $ws 12 4 $row11[3]
Set-TextValue  $ws 12 5 $row11[4]
Set-TextValue  $ws 12 6 $row11[5]

Write-Host "Done."
